# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    3  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    5  = @(1.445647641019636, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 3.005019366241741)
    6  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    7  = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    8  = @(0.2881169905109251, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1.012145535086602)
    9  = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 13.86384647080068, 37.47995879822157)
    10 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    11 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    12 = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
